# Generate Report for Handoff
# Rotates the tracked markdown file (new GUID) and its generated xliff
# handoff artifacts, refreshes the handoff timestamps, and resets the
# "latest target/handback" columns now that a new handoff cycle has begun.

$wb = $excel.ActiveWorkbook

$oldGuid = "77ce5d41-ba5b-4fc6-97bf-6b987d6f2067"
$newGuid = "d437ff78-8d10-4762-a667-2770c7afbdc1"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

$oldZhXlf = "$oldGuid.4ecf726fea8c58579321d739c736d4ed892de1b9.zh-cn.xlf"
$newZhXlf = "$newGuid.677caf15a6db5d497cbb2f770244abc034c5da93.zh-cn.xlf"

$oldDeXlf = "$oldGuid.4ecf726fea8c58579321d739c736d4ed892de1b9.de-de.xlf"
$newDeXlf = "$newGuid.677caf15a6db5d497cbb2f770244abc034c5da93.de-de.xlf"

$newHandoffDate = "2016-09-05 01:06:34"
$newZhHandoffDatetime = "2016-09-05 01:06:30"
$neverDate = "0001-01-01 00:00:00"

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = "e2e\$newMd"
$wsOverview.Range("G2").Value = $newHandoffDate

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newMd"
    }
}

# --- zh-cn sheet -------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDatetime
$wsZh.Range("K2").Value = $neverDate

foreach ($hl in @($wsZh.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMd
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsZh.Range("I2").Formula = "'"
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Formula = "'"
$wsZh.Range("J2").Style = "Normal"

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# --- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHandoffDate
$wsDe.Range("K2").Value = $neverDate

foreach ($hl in @($wsDe.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMd
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsDe.Range("I2").Formula = "'"
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Formula = "'"
$wsDe.Range("J2").Style = "Normal"

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
